# Add functionality to short custtypedetails data:
# append 5 new rows (Employee/EMP, EMPLOYEE/EMP, Distributor/DIST,
# Institution/INST, Ex. HQ/Ex.HQ) below the existing "custtypedetails"
# style list, give the first four new rows a thin box border, widen
# column B to fit, and leave the selection on the newly typed cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data rows -------------------------------------------------
$ws.Range("A128").Value = "Employee"
$ws.Range("B128").Value = "EMP"

$ws.Range("A129").Value = "EMPLOYEE"
$ws.Range("B129").Value = "EMP"

$ws.Range("A130").Value = "Distributor"
$ws.Range("B130").Value = "DIST"

$ws.Range("A131").Value = "Institution                                                 "
$ws.Range("B131").Value = "INST"

$ws.Range("A132").Value = "Ex. HQ"
$ws.Range("B132").Value = "Ex.HQ"

# --- formatting: thin box border around the first four new rows ----
$borderRange = $ws.Range("A128:B131")
$borderRange.Borders.LineStyle = 1
$borderRange.Borders.Weight = 2

# --- column sizing ---------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.63

# --- view / selection state -----------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 113
$win.ScrollColumn = 1
$ws.Range("C130").Select() | Out-Null
